# "b rith item import" - add a new leading "type" column to the
# items_import_form template sheet, shifting all existing header
# columns one to the right (A:T -> B:U).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; existing columns A:T shift to B:U,
# carrying their values, shared-string refs and header style along.
$ws.Columns.Item(1).Insert() | Out-Null

# Give the new A1 header cell the same style as its neighbours (copy the
# now-shifted original first header cell's formatting into it) and then
# set its text to the new field name.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null
$ws.Range("A1").Value = "type"

# Restore a plain single-cell selection on D5 (matches the saved view state).
$ws.Range("D5").Select() | Out-Null
